# Threat Alert Report update - 2026-01-25 09:07
# Refresh row 2 with the latest market scan, and append two new rows (3 & 4)
# with the newest fare-threat entries for SM-446 vs. Air Arabia Egypt E5-512.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Date column stores literal text like "13-MAR-26" (not real date
# serials), so pre-format column A as Text before writing the new date
# strings - otherwise Excel would auto-convert them into date serials.
$ws.Range("A2:A4").NumberFormat = "@"

# --- Update existing row 2 with refreshed figures ---
$ws.Range("A2").Value = "30-JAN-26"
$ws.Range("D2").Value = 612
$ws.Range("E2").Value = 895
$ws.Range("F2").Value = -283
$ws.Range("K2").Value = "SAR"

# --- Append row 3 ---
$ws.Range("A3").Value = "20-FEB-26"
$ws.Range("B3").Value = "SM-446"
$ws.Range("C3").Value = "Air Arabia Egypt E5-512"
$ws.Range("D3").Value = 513
$ws.Range("E3").Value = 883
$ws.Range("F3").Value = -370
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "LOW THREAT"
$ws.Range("K3").Value = "SAR"

# --- Append row 4 ---
$ws.Range("A4").Value = "27-FEB-26"
$ws.Range("B4").Value = "SM-446"
$ws.Range("C4").Value = "Air Arabia Egypt E5-512"
$ws.Range("D4").Value = 513
$ws.Range("E4").Value = 786
$ws.Range("F4").Value = -273
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "LOW THREAT"
$ws.Range("K4").Value = "SAR"

# --- Normalize formatting across the new/changed rows so they match the
#     existing table style (borders, centered alignment, fonts, fills) ---

# General data columns (A:I, K) take the plain bordered/centered style
# already used elsewhere in row 2 (e.g. B2), overriding the Text number
# format applied above now that the literal date strings are safely in.
$ws.Range("B2").Copy()
$ws.Range("A2:I4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K2:K4").PasteSpecial(-4122)   # xlPasteFormats

# IMPACT column keeps the bold "LOW THREAT" highlight style from J2.
$ws.Range("J2").Copy()
$ws.Range("J3:J4").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false
